# Add a new column X ("09-10-2020") to the COVID19 active-cases sheet,
# mirroring the existing daily-snapshot columns (D..W).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: a text label "09-10-2020" (not an actual date value), same
# as the other header cells in row 1 (e.g. W1 = "08-10-2020"). The leading
# apostrophe forces Excel to store it as literal text instead of parsing it
# into a date serial number.
$ws.Cells.Item(1, 24).Value = "'09-10-2020"

# Rows 2-36: new active-case counts per state/UT for 09-10-2020.
$data = @(
    @(2,  190),
    @(3,  48661),
    @(4,  2778),
    @(5,  30767),
    @(6,  11447),
    @(7,  1392),
    @(8,  27427),
    @(9,  105),
    @(10, 22232),
    @(11, 4716),
    @(12, 16465),
    @(13, 10867),
    @(14, 2943),
    @(15, 11482),
    @(16, 9272),
    @(17, 117162),
    @(18, 90664),
    @(19, 1299),
    @(20, 16788),
    @(21, 242438),
    @(22, 2877),
    @(23, 2369),
    @(24, 220),
    @(25, 1155),
    @(26, 26184),
    @(27, 4727),
    @(28, 10775),
    @(29, 21382),
    @(30, 545),
    @(31, 44437),
    @(32, 26374),
    @(33, 4197),
    @(34, 7849),
    @(35, 42552),
    @(36, 28854)
)

foreach ($pair in $data) {
    $row = $pair[0]
    $value = $pair[1]
    $ws.Cells.Item($row, 24).Value = $value
}
